$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a weekly price log for "Zapallo" (Camote / Paine varieties)
# at the Macroferia Regional de Talca. The update adds two new weekly
# observation rows (one for each variety) at two different insertion points,
# pushing the existing historical rows further down, and grows the sheet
# from A1:R262 to A1:R266.
# ---------------------------------------------------------------------------

function Set-Row($rowNum, $a,$b,$c,$d,$e,$f,$g,$h,$i,$j,$k,$l,$m,$n,$o,$p,$q,$r) {
    $ws.Range("A$rowNum").Value = $a
    $ws.Range("B$rowNum").Value = $b
    $ws.Range("C$rowNum").Value = $c
    $ws.Range("D$rowNum").Value = $d
    $ws.Range("E$rowNum").Value = $e
    $ws.Range("F$rowNum").Value = $f
    $ws.Range("G$rowNum").Value = $g
    $ws.Range("H$rowNum").Value = $h
    $ws.Range("I$rowNum").Value = $i
    $ws.Range("J$rowNum").Value = $j
    $ws.Range("K$rowNum").Value = $k
    $ws.Range("L$rowNum").Value = $l
    $ws.Range("M$rowNum").Value = $m
    $ws.Range("N$rowNum").Value = $n
    $ws.Range("O$rowNum").Value = $o
    $ws.Range("P$rowNum").Value = $p
    $ws.Range("Q$rowNum").Value = $q
    $ws.Range("R$rowNum").Value = $r
}

$unidad = "`$/kilo (volumen en unidades)"

# --- Insert the first pair of new rows at 226:227 (shifts old 226.. down by 2)
$ws.Rows("226:227").Insert()

Set-Row 226 5 "Macroferia Regional de Talca" "Maule" 44748 7 100112045 "Zapallo" "Camote" "1a (guarda)" 900 400 400 400 $unidad "Región del Maule" 400 1 "Hortaliza"
Set-Row 227 5 "Macroferia Regional de Talca" "Maule" 44748 7 100112045 "Zapallo" "Paine"  "1a (guarda)" 1500 150 150 150 $unidad "Región del Maule" 150 1 "Hortaliza"

# --- Insert the second pair of new rows at 255:256 (in the now-current row
#     numbering), shifting the remaining old rows (now at 255..262) down by 2
$ws.Rows("255:256").Insert()

Set-Row 255 5 "Macroferia Regional de Talca" "Maule" 44747 7 100112045 "Zapallo" "Camote" "1a (guarda)" 900 400 400 400 $unidad "Región del Maule" 400 1 "Hortaliza"
Set-Row 256 5 "Macroferia Regional de Talca" "Maule" 44747 7 100112045 "Zapallo" "Paine"  "1a (guarda)" 1500 150 150 150 $unidad "Región del Maule" 150 1 "Hortaliza"

Write-Output "Rows inserted; sheet now spans to row $($ws.Range("A1").CurrentRegion.Rows.Count)"
